$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.153.44'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '3.084.29'
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''522.85'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").Value = '''135.97'
$ws.Range("E6").Value = '  -4.32%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.083.73'
$ws.Range("E8").Value = '  -1.08%  '
$ws.Range("E9").Value = '  +2.42%  '
$ws.Range("D10").Value = '''7.34'
$ws.Range("E10").Value = '  +2.08%  '
$ws.Range("E11").Value = '  -2.46%  '
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = '3.625.68'
$ws.Range("E13").Value = '  -0.65%  '
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").Value = '''25.27'
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").Value = '''0.0000161'
$ws.Range("E16").Value = '  -2.67%  '
$ws.Range("D17").Value = '57.285.80'
$ws.Range("E17").Value = '  -1.63%  '
$ws.Range("D18").Value = '3.091.97'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").Value = '''5.87'
$ws.Range("E19").Value = '  -4.18%  '
$ws.Range("D20").Value = '''12.44'
$ws.Range("E20").Value = '  -2.70%  '
$ws.Range("D21").Value = '''7.84'
$ws.Range("E21").Value = '  -2.02%  '
$ws.Range("D22").Value = '''348.44'
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '''68.64'
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").Value = '''0.498'
$ws.Range("E25").Value = '  -3.44%  '
$ws.Range("E26").Value = '  -2.01%  '
$ws.Range("D27").Value = '''0.998'
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("D28").Value = '0.0₃0868'
$ws.Range("E28").Value = '  -6.84%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '''7.20'
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").Value = '''5.86'
$ws.Range("E32").Value = '  -8.32%  '
$ws.Range("D33").Value = '''20.85'
$ws.Range("E33").Value = '  -1.29%  '
$ws.Range("D34").Value = '''4.87'
$ws.Range("E34").Value = '  +4.72%  '
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = '''159.07'
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").Value = '''1.13'
$ws.Range("E36").Value = '  -4.94%  '
$ws.Range("D37").Value = '''6.02'
$ws.Range("E37").Value = '  -2.73%  '
$ws.Range("D38").Value = '''25.56'
$ws.Range("E38").Value = '  -3.03%  '
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("D40").Value = '''0.0658'
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("D41").Value = '''1.58'
$ws.Range("E41").Value = '  -2.96%  '
$ws.Range("D42").Value = '''4.03'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").Value = '''0.693'
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("D44").Value = '2.397.32'
$ws.Range("E44").Value = '  +5.51%  '
$ws.Range("D45").Value = '''36.66'
$ws.Range("E45").Value = '  -0.47%  '
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").Value = '3.131.97'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").Value = '''0.952'
$ws.Range("E49").Value = '  -4.71%  '
$ws.Range("D50").Value = '''5.96'
$ws.Range("E50").Value = '  -3.08%  '
$ws.Range("D51").Value = '''19.72'
$ws.Range("E51").Value = '  -4.46%  '
